$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws2 = $ws.Copy()
